$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster row: name + hyperlinked e-mail, matching the style of the
# existing rows (e.g. row 4 - Steven Koe).
$ws.Range("A5").Value = "Kiat Beng Goh"
$ws.Range("B5").Value = "kiatbeng.goh@digipen.edu"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:kiatbeng.goh@digipen.edu")

# Column A was manually widened (no longer auto "best fit").
$ws.Columns("A").ColumnWidth = 13.3

# Active selection moved to the newly-filled cell.
$ws.Range("B5").Select()
